$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New narrative text cells (order matters: this determines the shared-string
# table index assignment so it matches the target workbook exactly) ---
$ws.Range("E21").Value = "Esittely"
$ws.Range("E21").HorizontalAlignment = -4108   # xlCenter

$ws.Range("T22").Value = "Esittely ja vaatimusmäärittely"
$ws.Range("T22").HorizontalAlignment = -4108

$ws.Range("T21").Value = "Vaatimusmäärittely"
$ws.Range("T21").HorizontalAlignment = -4108

$ws.Range("O22").Value = "Ohjelmistokehitys, skaalaus, Retro, UML"

$ws.Range("T23").Value = "Vaatimusmäärittely, Retro"
$ws.Range("T23").HorizontalAlignment = -4108

$ws.Range("E22").Value = " Retro"
$ws.Range("E22").HorizontalAlignment = -4108

$ws.Range("J24").Value = " Retro"
$ws.Range("J24").HorizontalAlignment = -4108

$ws.Range("J23").Value = "Json ja .gitignoren ja solmuun mennen ohjelman korjailua"
$ws.Range("J23").HorizontalAlignment = -4108

# --- New logged hours / dates for the Retro entries ---
$ws.Range("C21").Value = 45356
$ws.Range("D21").Value = 1

$ws.Range("R21").Value = 45354
$ws.Range("S21").Value = 2

$ws.Range("C22").Value = 45358
$ws.Range("D22").Value = 4

$ws.Range("R22").Value = 45356
$ws.Range("S22").Value = 3

$ws.Range("N22").Value = 6

$ws.Range("R23").Value = 45358
$ws.Range("S23").Value = 5

$ws.Range("H24").Value = 45358
$ws.Range("I24").Value = 4

# --- View state ---
$ws.Range("J27").Select()
$ws.Columns.Item(15).ColumnWidth = 58.8
